# "update scripts wuth new tpm"
#
# The NATMI ligand/receptor pair table (Il13 -> Il2rg) was regenerated with
# updated TPM input data. The "MuSCs" cluster's underlying expression figures
# were recomputed as the "Resolving-Mac" cluster (the sending-cluster label in
# rows 6-9 changes from "MuSCs" to "Resolving-Mac"), and every derived
# specificity / weight column that depends on the ligand or receptor
# expression totals is refreshed to match the new TPM numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs (Il13) -> ECs (Il2rg)
$ws.Range("I2").Value = 0.9352601111131627
$ws.Range("J2").Value = 0.9352601111131628
$ws.Range("M2").Value = 35.32109533333334
$ws.Range("N2").Value = 105.963286
$ws.Range("O2").Value = 0.4123245624288747
$ws.Range("P2").Value = 0.4123245624288747
$ws.Range("Q2").Value = 11.33125463060067
$ws.Range("R2").Value = 101.981291675406
$ws.Range("S2").Value = 0.3856307160719155
$ws.Range("T2").Value = 0.3856307160719155

# Row 3: FAPs (Il13) -> FAPs (Il2rg)
$ws.Range("I3").Value = 0.9352601111131627
$ws.Range("J3").Value = 0.9352601111131628
$ws.Range("O3").Value = 0.01433703690686912
$ws.Range("P3").Value = 0.01433703690686912
$ws.Range("S3").Value = 0.01340885873055193
$ws.Range("T3").Value = 0.01340885873055193

# Row 4: FAPs (Il13) -> MuSCs (Il2rg)
$ws.Range("I4").Value = 0.9352601111131627
$ws.Range("J4").Value = 0.9352601111131628
$ws.Range("M4").Value = 2.583168333333334
$ws.Range("N4").Value = 7.749505
$ws.Range("O4").Value = 0.03015489023401347
$ws.Range("P4").Value = 0.03015489023401347
$ws.Range("Q4").Value = 0.8286984835116668
$ws.Range("R4").Value = 7.458286351605
$ws.Range("S4").Value = 0.02820266599086866
$ws.Range("T4").Value = 0.02820266599086866

# Row 5: FAPs (Il13) -> Resolving-Mac (Il2rg)
$ws.Range("I5").Value = 0.9352601111131627
$ws.Range("J5").Value = 0.9352601111131628
$ws.Range("M5").Value = 46.53090866666667
$ws.Range("N5").Value = 139.592726
$ws.Range("O5").Value = 0.5431835104302428
$ws.Range("P5").Value = 0.5431835104302427
$ws.Range("Q5").Value = 14.92744121662733
$ws.Range("R5").Value = 134.346970949646
$ws.Range("S5").Value = 0.5080178703198266
$ws.Range("T5").Value = 0.5080178703198266

# Row 6: sending cluster relabelled MuSCs -> Resolving-Mac, plus new TPM values; target ECs (Il2rg)
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.02220666666666667
$ws.Range("H6").Value = 0.06662
$ws.Range("I6").Value = 0.06473988888683736
$ws.Range("J6").Value = 0.06473988888683736
$ws.Range("M6").Value = 35.32109533333334
$ws.Range("N6").Value = 105.963286
$ws.Range("O6").Value = 0.4123245624288747
$ws.Range("P6").Value = 0.4123245624288747
$ws.Range("Q6").Value = 0.784363790368889
$ws.Range("R6").Value = 7.059274113320001
$ws.Range("S6").Value = 0.02669384635695918
$ws.Range("T6").Value = 0.02669384635695918

# Row 7: sending cluster relabelled MuSCs -> Resolving-Mac; target FAPs (Il2rg)
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("G7").Value = 0.02220666666666667
$ws.Range("H7").Value = 0.06662
$ws.Range("I7").Value = 0.06473988888683736
$ws.Range("J7").Value = 0.06473988888683736
$ws.Range("O7").Value = 0.01433703690686912
$ws.Range("P7").Value = 0.01433703690686912
$ws.Range("Q7").Value = 0.02727330272222222
$ws.Range("R7").Value = 0.2454597245
$ws.Range("S7").Value = 0.0009281781763171931
$ws.Range("T7").Value = 0.0009281781763171931

# Row 8: sending cluster relabelled MuSCs -> Resolving-Mac; target MuSCs (Il2rg)
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("G8").Value = 0.02220666666666667
$ws.Range("H8").Value = 0.06662
$ws.Range("I8").Value = 0.06473988888683736
$ws.Range("J8").Value = 0.06473988888683736
$ws.Range("M8").Value = 2.583168333333334
$ws.Range("N8").Value = 7.749505
$ws.Range("O8").Value = 0.03015489023401347
$ws.Range("P8").Value = 0.03015489023401347
$ws.Range("Q8").Value = 0.05736355812222223
$ws.Range("R8").Value = 0.5162720231
$ws.Range("S8").Value = 0.001952224243144809
$ws.Range("T8").Value = 0.001952224243144809

# Row 9: sending cluster relabelled MuSCs -> Resolving-Mac; target Resolving-Mac (Il2rg)
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 0.02220666666666667
$ws.Range("H9").Value = 0.06662
$ws.Range("I9").Value = 0.06473988888683736
$ws.Range("J9").Value = 0.06473988888683736
$ws.Range("M9").Value = 46.53090866666667
$ws.Range("N9").Value = 139.592726
$ws.Range("O9").Value = 0.5431835104302428
$ws.Range("P9").Value = 0.5431835104302427
$ws.Range("Q9").Value = 1.033296378457778
$ws.Range("R9").Value = 9.299667406119999
$ws.Range("S9").Value = 0.03516564011041618
$ws.Range("T9").Value = 0.03516564011041617
